$d = $word.ActiveDocument

$replacements = @(
    @{old = "50×60="; new = "95×95="},
    @{old = "36×74="; new = "50×72="},
    @{old = "28×38="; new = "60×45="},
    @{old = "24×13="; new = "89×47="},
    @{old = "62×63="; new = "69×16="},
    @{old = "39×36="; new = "32×14="},
    @{old = "41×87="; new = "30×51="},
    @{old = "92×87="; new = "82×49="},
    @{old = "70×81="; new = "41×50="},
    @{old = "94×50="; new = "48×66="},
    @{old = "44×54="; new = "52×15="},
    @{old = "31×97="; new = "53×42="},
    @{old = "34×63="; new = "11×15="},
    @{old = "32×47="; new = "84×71="},
    @{old = "39×56="; new = "90×74="},
    @{old = "98×63="; new = "35×37="},
    @{old = "47×45="; new = "81×73="},
    @{old = "23×95="; new = "35×39="},
    @{old = "51×32="; new = "71×83="},
    @{old = "89×43="; new = "71×91="},
    @{old = "95×25="; new = "98×83="},
    @{old = "20×13="; new = "65×14="},
    @{old = "58×49="; new = "33×34="},
    @{old = "61×21="; new = "40×69="},
    @{old = "50×21="; new = "31×80="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
